# Generate Report for Handback
# Refresh the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps for the last row (8a10b4e1-...7b883f13... handoff) on both the
# zh-cn and de-de language sheets, reflecting the latest report run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D14").Value = "2016-03-02 15:31:04"
$wsZhCn.Range("G14").Value = "2016-03-02 15:32:09"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D14").Value = "2016-03-02 15:31:31"
$wsDeDe.Range("G14").Value = "2016-03-02 15:32:28"
